# Daily attendance processing - 2026-01-07 18:44:38
# Swap the order of recorders in the "Recorded By" (column G) cells that
# contain both "dnasr281@gmail.com" and "System", so that "System" is
# listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
